$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.022.94'
$ws.Range('E2').Value = '  +5.33%  '
$ws.Range('D3').Value = '2.416.91'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').Value = '  +0.94%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.60'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.14'
$ws.Range('E6').Value = '  +6.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('E8').Value = '  +2.00%  '
$ws.Range('D9').Value = '2.448.77'
$ws.Range('E9').Value = '  +3.34%  '
$ws.Range('E10').Value = '  +5.93%  '
$ws.Range('E11').Value = '  +0.73%  '
$ws.Range('E12').Value = '  +3.04%  '
$ws.Range('E13').Value = '  +4.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.39'
$ws.Range('E14').Value = '  +6.90%  '
$ws.Range('E15').Value = '  +7.83%  '
$ws.Range('D16').Value = '2.858.20'
$ws.Range('E16').Value = '  +2.11%  '
$ws.Range('D17').Value = '62.865.00'
$ws.Range('E17').Value = '  +5.18%  '
$ws.Range('D18').Value = '2.452.00'
$ws.Range('E18').Value = '  +3.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.87'
$ws.Range('E19').Value = '  -2.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.99'
$ws.Range('E20').Value = '  +4.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '328.64'
$ws.Range('E21').Value = '  +2.34%  '
$ws.Range('E22').Value = '  +2.22%  '
$ws.Range('E23').Value = '  +14.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.55'
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '634.42'
$ws.Range('E26').Value = '  +13.77%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.10'
$ws.Range('E27').Value = '  +10.22%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.57'
$ws.Range('E28').Value = '  +5.56%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0988'
$ws.Range('E29').Value = '  +7.17%  '
$ws.Range('B30').Value = 'WrappedeETH'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D30').Value = '2.540.11'
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.24'
$ws.Range('E31').Value = '  +2.46%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.42'
$ws.Range('E32').Value = '  +8.98%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.139'
$ws.Range('E33').Value = '  +6.46%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.84'
$ws.Range('E34').Value = '  +3.41%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').Value = '  +4.45%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.995'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.76'
$ws.Range('E37').Value = '  +4.88%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.374'
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '153.01'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.42'
$ws.Range('E40').Value = '  +8.97%  '
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.69'
$ws.Range('E41').Value = '  +3.04%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.76'
$ws.Range('E42').Value = '  +13.45%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.77'
$ws.Range('E43').Value = '  +8.07%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0285'
$ws.Range('E45').Value = '  -4.46%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '145.19'
$ws.Range('E46').Value = '  +4.06%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.60'
$ws.Range('E47').Value = '  +2.53%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.48'
$ws.Range('E48').Value = '  +7.64%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.602'
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('B50').Value = 'Hedera'
$ws.Range('C50').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0517'
$ws.Range('E50').Value = '  +3.50%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '12.68'
$ws.Range('E51').Value = '  +8.57%  '
